$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.203.59"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "3.330.75"
$ws.Range("E3").Value = "  +7.82%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'252.90"
$ws.Range("E5").Value = "  +7.87%  "
$ws.Range("D6").Value = "'620.35"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").Value = "'0.383"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "3.319.30"
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "97.023.98"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("D15").Value = "'35.27"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").Value = "3.916.69"
$ws.Range("E16").Value = "  +6.88%  "
$ws.Range("D17").Value = "'5.51"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").Value = "3.333.46"
$ws.Range("E18").Value = "  +6.81%  "
$ws.Range("D19").Value = "'3.55"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").Value = "'14.86"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "'481.05"
$ws.Range("E21").Value = "  +8.82%  "
$ws.Range("E22").Value = "  +9.13%  "
$ws.Range("D23").Value = "'5.78"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'9.22"
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").Value = "'5.64"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("D26").Value = "'87.55"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").Value = "'12.02"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "3.502.15"
$ws.Range("E28").Value = "  +7.45%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").Value = "'0.239"
$ws.Range("E31").Value = "  -5.96%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").Value = "'9.17"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").Value = "'27.17"
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("D36").Value = "'7.38"
$ws.Range("E36").Value = "  -5.43%  "
$ws.Range("D37").Value = "'0.150"
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("D38").Value = "'506.42"
$ws.Range("E38").Value = "  +8.43%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").Value = "'0.448"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "'1.28"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'3.29"
$ws.Range("E43").Value = "  +5.44%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.794"
$ws.Range("E44").Value = "  +16.71%  "
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("D47").Value = "'161.24"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'1.91"
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("D49").Value = "'1.37"
$ws.Range("E49").Value = "  +6.28%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'45.43"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.49"
$ws.Range("E51").Value = "  +3.96%  "
